$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain text in this workbook (e.g. "51.960.91",
# "2.930.57") even though many values look numeric. Temporarily force the
# whole column to Text format so Excel's COM layer doesn't auto-convert
# numeric-looking strings (e.g. "357.52", "0.120") into real numbers (which
# would also silently drop significant trailing zeros). The format is
# restored to Normal afterwards so no visible formatting change remains.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "51.960.91"
$ws.Range("E2").Value = "  -0.88%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.930.57"
$ws.Range("E3").Value = "  +0.18%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 - BNB
$ws.Range("D5").Value = "357.52"
$ws.Range("E5").Value = "  +1.47%  "

# Row 6 - Solana
$ws.Range("D6").Value = "110.77"
$ws.Range("E6").Value = "  -1.86%  "

# Row 7 - XRP
$ws.Range("D7").Value = "0.569"
$ws.Range("E7").Value = "  +1.43%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.02%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.632"
$ws.Range("E9").Value = "  +0.65%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "39.55"
$ws.Range("E10").Value = "  -1.56%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +1.74%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.89%  "

# Row 13 - Chainlink
$ws.Range("D13").Value = "19.71"
$ws.Range("E13").Value = "  -2.30%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "7.93"
$ws.Range("E14").Value = "  +1.11%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "3.392.19"
$ws.Range("E15").Value = "  +0.22%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.932.07"
$ws.Range("E16").Value = "  -0.33%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  -0.68%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "51.975.05"
$ws.Range("E18").Value = "  -0.90%  "

# Row 19 - ImmutableX
$ws.Range("D19").Value = "3.38"
$ws.Range("E19").Value = "  +1.14%  "

# Row 20 - Uniswap
$ws.Range("E20").Value = "  -1.87%  "

# Row 21 - InternetComputer(DFINITY)
$ws.Range("D21").Value = "14.12"
$ws.Range("E21").Value = "  -2.71%  "

# Row 22 - ShibaInu
$ws.Range("D22").Value = "0.0₃0985"
$ws.Range("E22").Value = "  +0.18%  "

# Row 23 - Litecoin
$ws.Range("D23").Value = "71.12"
$ws.Range("E23").Value = "  -0.05%  "

# Row 24 - BitcoinCash
$ws.Range("D24").Value = "270.57"

# Row 25 - PancakeSwap
$ws.Range("D25").Value = "2.84"
$ws.Range("E25").Value = "  +1.21%  "

# Row 26 - Kaspa
$ws.Range("D26").Value = "0.187"
$ws.Range("E26").Value = "  +13.05%  "

# Row 27 - EthereumClassic
$ws.Range("D27").Value = "27.15"
$ws.Range("E27").Value = "  +0.74%  "

# Row 28 - Filecoin
$ws.Range("D28").Value = "7.55"
$ws.Range("E28").Value = "  +15.40%  "

# Row 29 - Dai
$ws.Range("E29").Value = "  +0.04%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "0.107"
$ws.Range("E30").Value = "  +13.82%  "

# Row 31 - Cosmos
$ws.Range("E31").Value = "  -0.10%  "

# Row 32 - InjectiveProtocol
$ws.Range("E32").Value = "  +0.01%  "

# Row 33 - Toncoin
$ws.Range("E33").Value = "  +1.76%  "

# Row 34 - RenderToken
$ws.Range("D34").Value = "6.07"
$ws.Range("E34").Value = "  -1.81%  "

# Row 35 - OKB
$ws.Range("D35").Value = "52.28"
$ws.Range("E35").Value = "  -1.95%  "

# Row 36 - VeChain
$ws.Range("D36").Value = "0.0444"
$ws.Range("E36").Value = "  -1.74%  "

# Row 37 - FirstDigitalUSD
$ws.Range("E37").Value = "  +0.05%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  -2.18%  "

# Row 39 - Celestia
$ws.Range("D39").Value = "18.52"
$ws.Range("E39").Value = "  -1.54%  "

# Row 40 - ARBITRUM
$ws.Range("E40").Value = "  -3.15%  "

# Row 41 - Stacks
$ws.Range("D41").Value = "2.75"
$ws.Range("E41").Value = "  +0.99%  "

# Row 42 - Stellar
$ws.Range("D42").Value = "0.120"
$ws.Range("E42").Value = "  +2.61%  "

# Row 43 - EnergySwap
$ws.Range("D43").Value = "23.11"

# Row 44 - Monero
$ws.Range("D44").Value = "119.36"
$ws.Range("E44").Value = "  -2.82%  "

# Row 45 - WEMIXToken
$ws.Range("E45").Value = "  -1.69%  "

# Row 46 - was ApeXProtocol, now NEARProtocol (swapped with row 47)
$ws.Range("B46").Value = "NEARProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D46").Value = "3.49"
$ws.Range("E46").Value = "  -2.19%  "

# Row 47 - was NEARProtocol, now ApeXProtocol (swapped with row 46)
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").Value = "2.49"
$ws.Range("E47").Value = "  -3.91%  "

# Row 48 - Maker
$ws.Range("D48").Value = "2.139.69"
$ws.Range("E48").Value = "  -3.53%  "

# Row 49 - TheGraph
$ws.Range("D49").Value = "0.252"
$ws.Range("E49").Value = "  -4.13%  "

# Row 50 - BEAM
$ws.Range("D50").Value = "0.0335"
$ws.Range("E50").Value = "  -0.96%  "

# Row 51 - FraxShare
$ws.Range("D51").Value = "9.18"
$ws.Range("E51").Value = "  -0.13%  "

# Restore the column's original (default) style now that all text values are
# safely in place.
$ws.Range("D2:D51").Style = "Normal"
